$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-03 Wednesday" "2025-09-04 Thursday"

Replace-Text "817×8=6536" "692×9=6228"
Replace-Text "779×9=7011" "630×5=3150"
Replace-Text "936×8=7488" "553×8=4424"
Replace-Text "220×6=1320" "980×3=2940"
Replace-Text "984×6=5904" "871×3=2613"
Replace-Text "564×2=1128" "948×2=1896"
Replace-Text "944×2=1888" "168×3=504"
Replace-Text "555×5=2775" "843×2=1686"
Replace-Text "357×6=2142" "675×5=3375"
Replace-Text "982×4=3928" "962×6=5772"
Replace-Text "595×6=3570" "984×2=1968"
Replace-Text "931×5=4655" "682×9=6138"
Replace-Text "285×7=1995" "347×7=2429"
Replace-Text "816×3=2448" "234×9=2106"
Replace-Text "755×5=3775" "455×9=4095"
Replace-Text "825×4=3300" "933×8=7464"
Replace-Text "494×5=2470" "565×5=2825"
Replace-Text "633×4=2532" "853×5=4265"
Replace-Text "553×3=1659" "611×5=3055"
Replace-Text "261×3=783" "782×7=5474"
Replace-Text "850×2=1700" "920×6=5520"
Replace-Text "951×8=7608" "601×8=4808"
Replace-Text "476×9=4284" "791×2=1582"
Replace-Text "492×9=4428" "844×4=3376"
Replace-Text "738×9=6642" "699×2=1398"
